$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.478.98"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "1.662.27"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("D4").Value = "'0.9993"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'234.96"
$ws.Range("E5").Value = "  +0.55%  "
$ws.Range("D7").Value = "'0.4620"
$ws.Range("E7").Value = "  -3.32%  "
$ws.Range("D8").Value = "'0.2567"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("D9").Value = "'0.06134"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("D10").Value = "1.658.94"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'0.06936"
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("D12").Value = "'14.59"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "'4.323"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").Value = "'75.02"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "'0.5688"
$ws.Range("E15").Value = "  -4.68%  "
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("D18").Value = "25.479.16"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "'0.000006675"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "'11.35"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "1.874.41"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("D22").Value = "'4.401"
$ws.Range("E22").Value = "  +1.65%  "
$ws.Range("D23").Value = "'8.629"
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "'5.205"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("D25").Value = "'134.14"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("D28").Value = "'1.706"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("D29").Value = "'103.60"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'3.935"
$ws.Range("E30").Value = "  +1.38%  "
$ws.Range("D31").Value = "'0.07697"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'3.585"
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("D33").Value = "'0.04330"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").Value = "'2.613"
$ws.Range("E34").Value = "  +1.66%  "
$ws.Range("D35").Value = "'0.5992"
$ws.Range("E35").Value = "  +0.85%  "
$ws.Range("D36").Value = "'0.9368"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "'0.9097"
$ws.Range("E37").Value = "  +5.24%  "
$ws.Range("B38").Value = "Quant"
$ws.Range("C38").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D38").Value = "'107.46"
$ws.Range("E38").Value = "  +8.77%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.398"
$ws.Range("E39").Value = "  -6.97%  "
$ws.Range("D40").Value = "'0.9993"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'0.01452"
$ws.Range("E41").Value = "  -3.92%  "
$ws.Range("D42").Value = "'1.814"
$ws.Range("E42").Value = "  +2.78%  "
$ws.Range("D43").Value = "'0.3698"
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").Value = "'4.986"
$ws.Range("E44").Value = "  +6.64%  "
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").Value = "'0.05259"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").Value = "'6.097"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").Value = "'30.37"
$ws.Range("E48").Value = "  +4.42%  "
$ws.Range("D49").Value = "'7.620"
$ws.Range("E49").Value = "  +7.28%  "
$ws.Range("D51").Value = "'0.9982"
$ws.Range("E51").Value = "  +0.11%  "
